$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Start from a clean slate for the region that is being restructured
#    (A1:G18 covers the old header/table rows; rows 19-28 keep their row
#    numbers but column B gets a new alignment style).
# ---------------------------------------------------------------------------
$ws.Range("A1:G18").Clear()

# ---------------------------------------------------------------------------
# 1. Build the two helper formats we need:
#    - "plain" : the existing family-3 (等线) font used for the epsilon-style
#                row labels (already present in the workbook as cell A3's
#                style before we cleared it) -- recreate it from scratch.
#    - "bold"  : the same family-3 font, bold -- used for every header /
#                section-title cell introduced by this change.
# ---------------------------------------------------------------------------

# Seed a throw-away cell with the original plain family-3 font so we can
# copy that exact font definition around (keeps family=3 instead of the
# family=2 the engine defaults to when you set Font.Name directly).
$ws.Range("Z1").Value = "seed"
$ws.Range("Z1").Font.Name = "等线"
$ws.Range("Z1").Font.Size = 11
$ws.Range("Z1").Font.Bold = $false
$plainSeed = $ws.Range("A6")
$plainSeed.Value = "seed2"
$plainSeed.Font.Name = "等线"

# Cells that use the plain (non-bold) family-3 font: epsilon-style row labels
$plainCells = "A6,A7,A8,A11,A12,A13"
$plainSeed.Copy()
$ws.Range($plainCells).PasteSpecial(-4122)

# Cells that use the bold family-3 font: every header / section-title cell
$boldCells = "A1:G1,A2:G2,A4,A10,A18"
$plainSeed.Copy()
$ws.Range($boldCells).PasteSpecial(-4122)
$ws.Range($boldCells).Font.Bold = $true

$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# ---------------------------------------------------------------------------
# 2. Row 1 - primary header
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "naïve_model"
$ws.Range("C1").Value = "LAT_model"
$ws.Range("D1").Value = "adv_training_model(eps=0.2,alpha = 0.5)"
$ws.Range("E1").Value = "adv_training_model(eps=0.1,alpha = 0.5)"
$ws.Range("F1").Value = "new_adversarial_training(alpha = 0.5,std = 0.2)"
$ws.Range("G1").Value = "new_adversarial_training(alpha = 0.5,std = 0.2)"

# ---------------------------------------------------------------------------
# 3. Row 2 - secondary header / note row
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "note"
$ws.Range("B2").Value = "baseline"
$ws.Range("C2").Value = "逐层对抗训练模型"
$ws.Range("D2").Value = "传统对抗训练模型(eps=0.2,alpha = 0.5)"
$ws.Range("E2").Value = "传统对抗训练模型(eps=0.1,alpha = 0.5)"
$ws.Range("F2").Value = "改进版对抗训练模型(alpha = 0.5,std = 0.2,,method = 'fgsm')"
$ws.Range("G2").Value = "改进版对抗训练模型(alpha = 0.5,std = 0.2,,method = 'tfgsm')"

# Row 3 intentionally left blank

# ---------------------------------------------------------------------------
# 4. fgsm block (rows 4-8)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "fgsm"

$ws.Range("A5").Value = "clean"
$ws.Range("B5").Value = 0.97
$ws.Range("C5").Value = 0.99
$ws.Range("D5").Value = 0.97
$ws.Range("E5").Value = 0.97
$ws.Range("F5").Value = 0.97
$ws.Range("G5").Value = 0.97

$ws.Range("A6").Value = "ε=0.1"
$ws.Range("B6").Value = 0.72
$ws.Range("C6").Value = 0.97
$ws.Range("D6").Value = 0.97
$ws.Range("E6").Value = 0.96
$ws.Range("F6").Value = 0.94
$ws.Range("G6").Value = 0.94

$ws.Range("A7").Value = "ε=0.2"
$ws.Range("B7").Value = 0.28000000000000003
$ws.Range("C7").Value = 0.92
$ws.Range("D7").Value = 0.97
$ws.Range("E7").Value = 0.92
$ws.Range("F7").Value = 0.89
$ws.Range("G7").Value = 0.88

$ws.Range("A8").Value = "ε=0.3"
$ws.Range("B8").Value = 0.04
$ws.Range("C8").Value = 0.73
$ws.Range("D8").Value = 0.95
$ws.Range("E8").Value = 0.78
$ws.Range("F8").Value = 0.77
$ws.Range("G8").Value = 0.72

# Row 9 intentionally left blank

# ---------------------------------------------------------------------------
# 5. ifgsm block (rows 10-13)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "ifgsm"

$ws.Range("A11").Value = "ε=0.1"
$ws.Range("B11").Value = 0.56000000000000005
$ws.Range("C11").Value = 0.97
$ws.Range("D11").Value = 0.97
$ws.Range("E11").Value = 0.96
$ws.Range("F11").Value = 0.94
$ws.Range("G11").Value = 0.94

$ws.Range("A12").Value = "ε=0.15"
$ws.Range("B12").Value = 0.15
$ws.Range("C12").Value = 0.95
$ws.Range("D12").Value = 0.96
$ws.Range("E12").Value = 0.94
$ws.Range("F12").Value = 0.92
$ws.Range("G12").Value = 0.91

$ws.Range("A13").Value = "ε=0.2"
$ws.Range("B13").Value = 0.01
$ws.Range("C13").Value = 0.91
$ws.Range("D13").Value = 0.95
$ws.Range("E13").Value = 0.92
$ws.Range("F13").Value = 0.88
$ws.Range("G13").Value = 0.86

# Rows 14-17 intentionally left blank

# ---------------------------------------------------------------------------
# 6. hyperparameters block (rows 18-28) - row numbers unchanged from before
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "hyperparameters"

$ws.Range("A19").Value = "alpha"
$ws.Range("B19").Value = 0.6

$ws.Range("A20").Value = "batch_norm"
$ws.Range("B20").Value = $true

$ws.Range("A21").Value = "batch_size"
$ws.Range("B21").Value = 64

$ws.Range("A22").Value = "dataset"
$ws.Range("B22").Value = "MNIST"

$ws.Range("A23").Value = "dropout"
$ws.Range("B23").Value = $true

$ws.Range("A24").Value = "epoch"
$ws.Range("B24").Value = 2

$ws.Range("A25").Value = "epsilon"
$ws.Range("B25").Value = 0.6

$ws.Range("A26").Value = "input_ch"
$ws.Range("B26").Value = 3

$ws.Range("A27").Value = "lr"
$ws.Range("B27").Value = 0.0002

$ws.Range("A28").Value = "pro_num"
$ws.Range("B28").Value = 8

# Right-align the hyperparameter values column
$ws.Range("B19:B28").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 7. Cosmetics: column widths + selection, matching the post-edit workbook
# ---------------------------------------------------------------------------
$ws.Columns("A:G").AutoFit()

$ws.Range("F27").Select()

Write-Output "edit applied"
